$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimals (e.g. "234.11") that Excel's
# smart-entry would otherwise coerce to a Number; force those specific
# cells to Text, write the value, then restore the default style so no
# cell formatting changes are left behind.

$ws.Range("D2").Value = "37.253.50"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.068.81"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0763"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "2.374.80"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "2.069.83"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "37.224.85"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "0.0₃0812"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "1.480.48"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0935"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0211"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.27%  "
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "2.263.67"
$ws.Range("E51").Value = "  -0.48%  "
